# Refresh crypto price/volume table (includes a WrappedEther/ShibaInu row swap)
# to match the latest feed pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new "Price" values look like plain numbers (e.g. "586.65", "0.650").
# The Value setter would auto-convert those to floating point, dropping
# trailing zeros / losing precision and changing the cell from text to a
# number. Force text via NumberFormat "@" first, then restore the default
# "Normal" style afterwards so no stray number format sticks to the cell.

$ws.Range("D2").Value = "66.388.33"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "3.325.48"
$ws.Range("E3").Value = "  +0.01%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.65"
$ws.Range("E5").Value = "  +2.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.73"
$ws.Range("E6").Value = "  +0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("E7").Value = "  +4.89%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "3.325.03"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.81"
$ws.Range("E11").Value = "  +2.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.402"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "3.907.16"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("E14").Value = "  -2.29%  "

$ws.Range("D15").Value = "66.414.26"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.64"
$ws.Range("E16").Value = "  -0.22%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.381.50"
$ws.Range("E17").Value = "  +1.36%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000164"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "425.39"
$ws.Range("E19").Value = "  -2.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.53"
$ws.Range("E20").Value = "  -2.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -2.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.38"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.91"
$ws.Range("E23").Value = "  -1.73%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").Value = "3.474.57"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.514"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.206"
$ws.Range("E28").Value = "  +6.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000115"
$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.06"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  -1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.39"
$ws.Range("E33").Value = "  -1.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.19"
$ws.Range("E35").Value = "  -0.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.64"
$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("E37").Value = "  -2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.10"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("E39").Value = "  -2.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  +0.83%  "

$ws.Range("D41").Value = "2.863.89"
$ws.Range("E41").Value = "  +1.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.47"
$ws.Range("E42").Value = "  -4.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.34"
$ws.Range("E43").Value = "  -1.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.759"
$ws.Range("E44").Value = "  -4.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.79"
$ws.Range("E45").Value = "  -1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0663"
$ws.Range("E46").Value = "  -0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.99"
$ws.Range("E47").Value = "  -2.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -0.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.17"
$ws.Range("E49").Value = "  -3.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "311.99"
$ws.Range("E50").Value = "  -3.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0272"
$ws.Range("E51").Value = "  +0.60%  "

# Restore the default (no explicit number format) style on the cells above
# that we temporarily forced to text.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
